$d = $word.ActiveDocument

# Locate the paragraph containing the first question text and the paragraph
# containing the tail of the "Quem sera" question, then expand each Find hit
# to its full enclosing paragraph so the replacement range spans exactly the
# two original list-item paragraphs that need restructuring.
$r1 = $d.Content
$null = $r1.Find.Execute("Como simular o fluxo do tr", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r1.Expand(4)

$r2 = $d.Content
$null = $r2.Find.Execute("Alguma empresa por tr", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r2.Expand(4)

$target = $d.Range($r1.Start, $r2.End)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:ind w:left="284"/><w:jc w:val="both"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:ind w:left="284" w:hanging="284"/><w:jc w:val="both"/></w:pPr><w:r><w:t>Como simular o fluxo do tráfego de uma cidade? Como obter os dados reais a fim de simulá-los? Há algum padrão de comportamento?</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="PargrafodaLista"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:ind w:left="284"/><w:jc w:val="both"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:ind w:left="284" w:hanging="284"/><w:jc w:val="both"/></w:pPr><w:r><w:t>Quem será o responsável pela utilização do software?</w:t></w:r><w:r><w:t xml:space="preserve"> Alguma empresa por trás disso?</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:ind w:left="284"/><w:jc w:val="both"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:ind w:left="284" w:hanging="284"/><w:jc w:val="both"/></w:pPr><w:r><w:t>Há algum padrão de E/S?</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target.InsertXML($xml)

# Best-effort: mark the "Default Paragraph Font" character style as semi-hidden
# (w:semiHidden) to mirror the styles.xml update in the source change.
$style = $d.Styles.Item("Fontepargpadro")
try {
    $style.Hidden = $true
} catch {
}

Write-Output ("Paragraphs.Count=" + $d.Paragraphs.Count)
